# Title Page.docx edits
#
# 1) Affiliations paragraph: split the Dept. of Pediatrics run right after
#    "University of T" and drop Word's "_GoBack" bookmark at that point
#    (simulates the author clicking/editing there last).
# 2) Keywords line: insert " syndrome" after "Cushing's" (split into 3 runs).
# 3) Word Count line: change "4511" -> "4512" (split into 3 runs).
# 4) The old "_GoBack" bookmark (previously sitting alone in the final
#    empty paragraph) moves with the edit above, so that trailing
#    paragraph reverts to being a plain empty paragraph.

$d = $word.ActiveDocument

function Split-RunAt {
    # Force a run boundary at an absolute document character offset by
    # dropping a throwaway bookmark there and immediately deleting it.
    # Word keeps the run split even after the bookmark itself is gone.
    param($doc, [int]$pos, [string]$tmpName)

    $ptRange = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add($tmpName, $ptRange)
    $doc.Bookmarks($tmpName).Delete()
}

# ---------------------------------------------------------------------
# Edit 2: Keywords -> insert " syndrome" after "Cushing's"
# ---------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Cushing’s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    $find = $d.Content
    $found = $find.Find.Execute("Cushing's", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

$insertPt = $find.Duplicate
$insertPt.Collapse(0)
$beforePos = $insertPt.Start
$insertPt.InsertAfter(" syndrome")
$afterPos = $insertPt.End

Split-RunAt $d $beforePos "TmpSplitA"
Split-RunAt $d $afterPos "TmpSplitB"

# ---------------------------------------------------------------------
# Edit 3: Word Count 4511 -> 4512
# ---------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("4511", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$digitStart = $find2.Start + 3   # position right before the final "1"
$digitEnd = $find2.Start + 4     # position right after the final "1"

$digitRange = $d.Range($digitStart, $digitEnd)
$digitRange.Text = "2"

Split-RunAt $d $digitStart "TmpSplitC"
Split-RunAt $d ($digitStart + 1) "TmpSplitD"

# ---------------------------------------------------------------------
# Edit 1: Affiliations - split "...University of T|ennessee..." and drop
# the real "_GoBack" bookmark at the split point (this also removes it
# from wherever it used to be, i.e. the trailing empty paragraph).
# ---------------------------------------------------------------------
$find3 = $d.Content
$found3 = $find3.Find.Execute(" Department of Pediatrics, University of T", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$gobackPt = $find3.Duplicate
$gobackPt.Collapse(0)
$d.Bookmarks.Add("_GoBack", $gobackPt)

Write-Output "done"
